$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1234
$ws.Range("F7").Value = 1415
$ws.Range("F8").Value = 76
$ws.Range("F9").Value = 20
$ws.Range("F10").Value = 661
$ws.Range("F11").Value = 141
$ws.Range("F12").Value = 136
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 499
$ws.Range("F20").Value = 110
$ws.Range("F21").Value = 726
$ws.Range("F22").Value = 2552
$ws.Range("F26").Value = 295
$ws.Range("F28").Value = 9
$ws.Range("F29").Value = 111
$ws.Range("F30").Value = 568
$ws.Range("F32").Value = 45
$ws.Range("F33").Value = 101
$ws.Range("F38").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 25
$ws.Range("F5").Value = 625
$ws.Range("F6").Value = 625
$ws.Range("F16").Value = 487
$ws.Range("F19").Value = 943
$ws.Range("F26").Value = 244
$ws.Range("F29").Value = 1
$ws.Range("F31").Value = 188
$ws.Range("F33").Value = 19

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1756
$ws.Range("F5").Value = 2297
$ws.Range("F6").Value = 917
$ws.Range("F9").Value = 1150
$ws.Range("F10").Value = 274

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1756
$ws.Range("F4").Value = 2297
$ws.Range("F7").Value = 25
$ws.Range("F8").Value = 917
$ws.Range("F9").Value = 1150
$ws.Range("F10").Value = 274
$ws.Range("F12").Value = 1234
$ws.Range("F15").Value = 1415
$ws.Range("F16").Value = 625
$ws.Range("F17").Value = 76
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 661
$ws.Range("F20").Value = 141
$ws.Range("F22").Value = 136
$ws.Range("F23").Value = 499
$ws.Range("F26").Value = 110
$ws.Range("F27").Value = 726
$ws.Range("F28").Value = 2552
$ws.Range("F31").Value = 295
$ws.Range("F34").Value = 111
$ws.Range("F36").Value = 568
$ws.Range("F38").Value = 487
$ws.Range("F40").Value = 45
$ws.Range("F41").Value = 101
$ws.Range("F45").Value = 244
$ws.Range("F47").Value = 188
$ws.Range("F50").Value = 27

